$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove obsolete columns J and O (now merged/removed), update all other values
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 2
$ws.Range("D2").Value = 49624
$ws.Range("E2").Value = 1433
$ws.Range("F2").Value = 1433
$ws.Range("G2").Value = 1486
$ws.Range("H2").Value = 1113
$ws.Range("I2").Value = 1113
$ws.Range("K2").Value = 29201
$ws.Range("L2").Value = 12385
$ws.Range("M2").Value = 16815
$ws.Range("N2").Value = 16815
$ws.Range("P2").Value = 770
$ws.Range("Q2").Value = 3684
$ws.Range("R2").Value = -1808
$ws.Range("S2").Value = -1276
$ws.Range("T2").Value = 1875
$ws.Range("U2").Value = 1810
$ws.Range("V2").Value = 5001
$ws.Range("W2").Value = 2.89
$ws.Range("X2").Value = 2.24
$ws.Range("Y2").Value = 6.75
$ws.Range("Z2").Value = 3.84
$ws.Range("AA2").Value = 73.66
$ws.Range("AB2").Value = 2084.7
$ws.Range("AC2").Value = 1445
$ws.Range("AD2").Value = 17.75
$ws.Range("AE2").Value = 21838
$ws.Range("AF2").Value = 1.17
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 2.34
$ws.Range("AI2").Value = 41.52
$ws.Range("AJ2").Value = 77000000

# Row 3
$ws.Range("D3").Value = 62731
$ws.Range("E3").Value = 2258
$ws.Range("F3").Value = 2258
$ws.Range("G3").Value = 2214
$ws.Range("H3").Value = 1662
$ws.Range("I3").Value = 1642
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 47355
$ws.Range("L3").Value = 25765
$ws.Range("M3").Value = 21590
$ws.Range("N3").Value = 17926
$ws.Range("O3").Value = 3664
$ws.Range("P3").Value = 770
$ws.Range("Q3").Value = 4440
$ws.Range("R3").Value = -8857
$ws.Range("S3").Value = 4026
$ws.Range("T3").Value = 2552
$ws.Range("U3").Value = 1888
$ws.Range("V3").Value = 11969
$ws.Range("W3").Value = 3.6
$ws.Range("X3").Value = 2.65
$ws.Range("Y3").Value = 9.449999999999999
$ws.Range("Z3").Value = 4.34
$ws.Range("AA3").Value = 119.33
$ws.Range("AB3").Value = 2228.83
$ws.Range("AC3").Value = 2133
$ws.Range("AD3").Value = 25.23
$ws.Range("AE3").Value = 23281
$ws.Range("AF3").Value = 2.31
$ws.Range("AG3").Value = 850
$ws.Range("AH3").Value = 1.58
$ws.Range("AI3").Value = 39.86
$ws.Range("AJ3").Value = 77000000

# Row 4
$ws.Range("D4").Value = 74020
$ws.Range("E4").Value = 2181
$ws.Range("F4").Value = 2181
$ws.Range("G4").Value = 3646
$ws.Range("H4").Value = 2735
$ws.Range("I4").Value = 2743
$ws.Range("J4").Value = -7
$ws.Range("K4").Value = 47214
$ws.Range("L4").Value = 23638
$ws.Range("M4").Value = 23575
$ws.Range("N4").Value = 19898
$ws.Range("O4").Value = 3677
$ws.Range("P4").Value = 770
$ws.Range("Q4").Value = 4345
$ws.Range("R4").Value = 78
$ws.Range("S4").Value = -4635
$ws.Range("T4").Value = 3110
$ws.Range("U4").Value = 1235
$ws.Range("V4").Value = 8190
$ws.Range("W4").Value = 2.95
$ws.Range("X4").Value = 3.69
$ws.Range("Y4").Value = 14.5
$ws.Range("Z4").Value = 5.79
$ws.Range("AA4").Value = 100.27
$ws.Range("AB4").Value = 2484.48
$ws.Range("AC4").Value = 3562
$ws.Range("AD4").Value = 13.36
$ws.Range("AE4").Value = 25842
$ws.Range("AF4").Value = 1.84
$ws.Range("AG4").Value = 1100
$ws.Range("AH4").Value = 2.31
$ws.Range("AI4").Value = 30.88
$ws.Range("AJ4").Value = 77000000

# Row 5
$ws.Range("D5").Value = 82666
$ws.Range("E5").Value = 1657
$ws.Range("F5").Value = 1657
$ws.Range("G5").Value = 1892
$ws.Range("H5").Value = 1151
$ws.Range("I5").Value = 1182
$ws.Range("J5").Value = -31
$ws.Range("K5").Value = 50923
$ws.Range("L5").Value = 26989
$ws.Range("M5").Value = 23933
$ws.Range("N5").Value = 20304
$ws.Range("O5").Value = 3630
$ws.Range("P5").Value = 770
$ws.Range("Q5").Value = 4422
$ws.Range("R5").Value = -4233
$ws.Range("S5").Value = 282
$ws.Range("T5").Value = 3809
$ws.Range("U5").Value = 614
$ws.Range("V5").Value = 9867
$ws.Range("W5").Value = 2
$ws.Range("X5").Value = 1.39
$ws.Range("Y5").Value = 5.88
$ws.Range("Z5").Value = 2.35
$ws.Range("AA5").Value = 112.77
$ws.Range("AB5").Value = 2541.87
$ws.Range("AC5").Value = 1535
$ws.Range("AD5").Value = 26.26
$ws.Range("AE5").Value = 26369
$ws.Range("AF5").Value = 1.53
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 1.49
$ws.Range("AI5").Value = 39.1
$ws.Range("AJ5").Value = 77000000

# Row 6
$ws.Range("D6").Value = 86916
$ws.Range("E6").Value = 1803
$ws.Range("F6").Value = 1803
$ws.Range("G6").Value = 1795
$ws.Range("H6").Value = 1323
$ws.Range("I6").Value = 1206
$ws.Range("K6").Value = 50989
$ws.Range("L6").Value = 26405
$ws.Range("M6").Value = 24583
$ws.Range("N6").Value = 20879
$ws.Range("P6").Value = 770
$ws.Range("Q6").Value = 3248
$ws.Range("R6").Value = -3849
$ws.Range("S6").Value = 311
$ws.Range("T6").Value = 2876
$ws.Range("U6").Value = 372
$ws.Range("V6").Value = 10711
$ws.Range("W6").Value = 2.07
$ws.Range("X6").Value = 1.52
$ws.Range("Y6").Value = 5.86
$ws.Range("Z6").Value = 2.6
$ws.Range("AA6").Value = 107.41
$ws.Range("AB6").Value = 2617.34
$ws.Range("AC6").Value = 1566
$ws.Range("AD6").Value = 25.86
$ws.Range("AE6").Value = 27116
$ws.Range("AF6").Value = 1.49
$ws.Range("AG6").Value = 650
$ws.Range("AH6").Value = 1.6
$ws.Range("AI6").Value = 41.5
$ws.Range("AJ6").Value = 77000000

# Row 7
$ws.Range("D7").Value = 90238
$ws.Range("E7").Value = 2247
$ws.Range("G7").Value = 2130
$ws.Range("H7").Value = 1578
$ws.Range("I7").Value = 1494
$ws.Range("K7").Value = 64366
$ws.Range("L7").Value = 38612
$ws.Range("M7").Value = 25754
$ws.Range("N7").Value = 21977
$ws.Range("P7").Value = 770
$ws.Range("Q7").Value = 6403
$ws.Range("R7").Value = -5534
$ws.Range("S7").Value = -770
$ws.Range("T7").Value = 2179
$ws.Range("U7").Value = 3771
$ws.Range("W7").Value = 2.49
$ws.Range("X7").Value = 1.75
$ws.Range("Y7").Value = 6.97
$ws.Range("Z7").Value = 2.74
$ws.Range("AA7").Value = 149.93
$ws.Range("AC7").Value = 1940
$ws.Range("AD7").Value = 20.44
$ws.Range("AE7").Value = 28542
$ws.Range("AF7").Value = 1.39
$ws.Range("AG7").Value = 730
$ws.Range("AH7").Value = 1.84
$ws.Range("AI7").Value = 37.63

# Row 8
$ws.Range("D8").Value = 94615
$ws.Range("E8").Value = 2599
$ws.Range("G8").Value = 2506
$ws.Range("H8").Value = 1886
$ws.Range("I8").Value = 1811
$ws.Range("K8").Value = 64623
$ws.Range("L8").Value = 37531
$ws.Range("M8").Value = 27092
$ws.Range("N8").Value = 23244
$ws.Range("P8").Value = 770
$ws.Range("Q8").Value = 5670
$ws.Range("R8").Value = -2480
$ws.Range("S8").Value = -2381
$ws.Range("T8").Value = 2027
$ws.Range("U8").Value = 2777
$ws.Range("W8").Value = 2.75
$ws.Range("X8").Value = 1.99
$ws.Range("Y8").Value = 8.01
$ws.Range("Z8").Value = 2.92
$ws.Range("AA8").Value = 138.53
$ws.Range("AC8").Value = 2352
$ws.Range("AD8").Value = 16.85
$ws.Range("AE8").Value = 30186
$ws.Range("AF8").Value = 1.31
$ws.Range("AG8").Value = 780
$ws.Range("AH8").Value = 1.97
$ws.Range("AI8").Value = 33.16

# Row 9
$ws.Range("D9").Value = 99190
$ws.Range("E9").Value = 2944
$ws.Range("G9").Value = 2893
$ws.Range("H9").Value = 2163
$ws.Range("I9").Value = 2081
$ws.Range("K9").Value = 65834
$ws.Range("L9").Value = 37289
$ws.Range("M9").Value = 28545
$ws.Range("N9").Value = 24624
$ws.Range("P9").Value = 770
$ws.Range("Q9").Value = 5667
$ws.Range("R9").Value = -2980
$ws.Range("S9").Value = -1721
$ws.Range("T9").Value = 2071
$ws.Range("U9").Value = 2987
$ws.Range("W9").Value = 2.97
$ws.Range("X9").Value = 2.18
$ws.Range("Y9").Value = 8.699999999999999
$ws.Range("Z9").Value = 3.32
$ws.Range("AA9").Value = 130.63
$ws.Range("AC9").Value = 2703
$ws.Range("AD9").Value = 14.67
$ws.Range("AE9").Value = 31979
$ws.Range("AF9").Value = 1.24
$ws.Range("AG9").Value = 830
$ws.Range("AH9").Value = 2.09
$ws.Range("AI9").Value = 30.71
